$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Enterprises density (per 1000 people)
$ws.Range("B11").Formula = "'35.88"
$ws.Range("C11").Formula = "'8.79"
$ws.Range("D11").Formula = "'44.66"

# Row 12: Employment (% of total)
$ws.Range("B12").Formula = "'26.98"
$ws.Range("C12").Formula = "'60.15"
$ws.Range("D12").Formula = "'87.13"

# Row 14: Enterprises (% of total)
$ws.Range("B14").Formula = "'80.18"
$ws.Range("C14").Formula = "'19.63"
$ws.Range("D14").Formula = "'99.81"
